$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.701.01"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.274.74"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  +3.84%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.63%  "
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.906"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.60%  "
$ws.Range("D16").Value = "2.615.89"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "2.272.45"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "43.631.28"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.11%  "
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("E36").Value = "  +13.73%  "
$ws.Range("E37").Value = "  +9.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.240"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.56%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.87%  "
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.682"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +22.57%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "74.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +39.74%  "
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.48%  "
